# Rename the "wt" and "dcin5" sheets to reflect that they now hold
# log2 expression data, and switch the active/selected tab from
# "optimization_parameters" to the renamed "wt" sheet
# ("wt_log2_expression").

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("wt").Name = "wt_log2_expression"
$wb.Worksheets.Item("dcin5").Name = "dcin5_log2_expression"

# Activating this sheet moves tabSelected from the previously active
# sheet (optimization_parameters) to this one, and updates the
# workbook's activeTab accordingly.
$wb.Worksheets.Item("wt_log2_expression").Activate()
